$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the looting item name in B3 from "rock" to "scrap"
$ws.Range("B3").Value = "scrap"

# Change the associated lootingId in D3 from 1 to 5001
$ws.Range("D3").Value = 5001

# Reflect the active selection left on D3 after the edit
$ws.Range("D3").Select()
